$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two rows for the new 2021-10-04 ("semana") records, shifting old row 13+ down.
$ws.Range("A13:A14").EntireRow.Insert()

# row 13
$ws.Cells.Item(13,1).Value = 9
$ws.Cells.Item(13,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(13,3).Value = 'Metropolitana'
$ws.Cells.Item(13,4).Value = 44473
$ws.Cells.Item(13,5).Value = 13
$ws.Cells.Item(13,6).Value = 'Fruta'
$ws.Cells.Item(13,7).Value = 100108
$ws.Cells.Item(13,8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(13,9).Value = 100108003
$ws.Cells.Item(13,10).Value = 'Maracuyá'
$ws.Cells.Item(13,11).Value = 'Sin especificar'
$ws.Cells.Item(13,12).Value = 'Primera'
$ws.Cells.Item(13,13).Value = 25
$ws.Cells.Item(13,14).Value = 72000
$ws.Cells.Item(13,15).Value = 72000
$ws.Cells.Item(13,16).Value = 72000
$ws.Cells.Item(13,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(13,18).Value = 'Perú'
$ws.Cells.Item(13,19).Value = 4000
$ws.Cells.Item(13,20).Value = 18

# row 14
$ws.Cells.Item(14,1).Value = 9
$ws.Cells.Item(14,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(14,3).Value = 'Metropolitana'
$ws.Cells.Item(14,4).Value = 44473
$ws.Cells.Item(14,5).Value = 13
$ws.Cells.Item(14,6).Value = 'Fruta'
$ws.Cells.Item(14,7).Value = 100108
$ws.Cells.Item(14,8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(14,9).Value = 100108003
$ws.Cells.Item(14,10).Value = 'Maracuyá'
$ws.Cells.Item(14,11).Value = 'Sin especificar'
$ws.Cells.Item(14,12).Value = 'Segunda'
$ws.Cells.Item(14,13).Value = 5
$ws.Cells.Item(14,14).Value = 67000
$ws.Cells.Item(14,15).Value = 67000
$ws.Cells.Item(14,16).Value = 67000
$ws.Cells.Item(14,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(14,18).Value = 'Perú'
$ws.Cells.Item(14,19).Value = 3722
$ws.Cells.Item(14,20).Value = 18

# Insert two more rows for the new 2021-09-27 records (post-shift position 30), shifting the rest down.
$ws.Range("A30:A31").EntireRow.Insert()

# row 30
$ws.Cells.Item(30,1).Value = 9
$ws.Cells.Item(30,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(30,3).Value = 'Metropolitana'
$ws.Cells.Item(30,4).Value = 44466
$ws.Cells.Item(30,5).Value = 13
$ws.Cells.Item(30,6).Value = 'Fruta'
$ws.Cells.Item(30,7).Value = 100108
$ws.Cells.Item(30,8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(30,9).Value = 100108003
$ws.Cells.Item(30,10).Value = 'Maracuyá'
$ws.Cells.Item(30,11).Value = 'Sin especificar'
$ws.Cells.Item(30,12).Value = 'Primera'
$ws.Cells.Item(30,13).Value = 15
$ws.Cells.Item(30,14).Value = 60000
$ws.Cells.Item(30,15).Value = 60000
$ws.Cells.Item(30,16).Value = 60000
$ws.Cells.Item(30,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(30,18).Value = 'Perú'
$ws.Cells.Item(30,19).Value = 3333
$ws.Cells.Item(30,20).Value = 18

# row 31
$ws.Cells.Item(31,1).Value = 9
$ws.Cells.Item(31,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(31,3).Value = 'Metropolitana'
$ws.Cells.Item(31,4).Value = 44466
$ws.Cells.Item(31,5).Value = 13
$ws.Cells.Item(31,6).Value = 'Fruta'
$ws.Cells.Item(31,7).Value = 100108
$ws.Cells.Item(31,8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(31,9).Value = 100108003
$ws.Cells.Item(31,10).Value = 'Maracuyá'
$ws.Cells.Item(31,11).Value = 'Sin especificar'
$ws.Cells.Item(31,12).Value = 'Segunda'
$ws.Cells.Item(31,13).Value = 10
$ws.Cells.Item(31,14).Value = 58000
$ws.Cells.Item(31,15).Value = 58000
$ws.Cells.Item(31,16).Value = 58000
$ws.Cells.Item(31,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(31,18).Value = 'Perú'
$ws.Cells.Item(31,19).Value = 3222
$ws.Cells.Item(31,20).Value = 18

